$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6952.5
$ws.Range("I43").Value = 4400.273
$ws.Range("J43").Value = 9112.077
$ws.Range("K43").Value = 4400.273
$ws.Range("L43").Value = 9112.077
$ws.Range("M43").Value = -4331.273
$ws.Range("N43").Value = -9250.077
$ws.Range("H58").Value = 2648.6428
$ws.Range("I58").Value = 864.0
$ws.Range("J58").Value = 3987.125
$ws.Range("K58").Value = 2592.0
$ws.Range("L58").Value = 11961.375
$ws.Range("M58").Value = -2442.0
$ws.Range("N58").Value = -12261.375
$ws.Range("H74").Value = 4800.3
$ws.Range("I74").Value = 4800.3
$ws.Range("K74").Value = 4800.3
$ws.Range("M74").Value = -3864.3
$ws.Range("H77").Value = 4800.3
$ws.Range("I77").Value = 4800.3
$ws.Range("K77").Value = 24001.5
$ws.Range("M77").Value = -19321.5
$ws.Range("H87").Value = 0.0
$ws.Range("J87").Value = 0.0
$ws.Range("L87").ClearContents()
$ws.Range("N87").Value = 0.0
$ws.Range("H90").Value = 0.0
$ws.Range("J90").Value = 0.0
$ws.Range("L90").ClearContents()
$ws.Range("N90").Value = 0.0
$ws.Range("H113").Value = 10183.0
$ws.Range("I113").Value = 21599.5
$ws.Range("K113").Value = 21599.5
$ws.Range("M113").Value = -18345.5
$ws.Range("H132").Value = 1964563.8
$ws.Range("I132").Value = 3912.205
$ws.Range("K132").Value = 11736.615
$ws.Range("M132").Value = -9206.615
$ws.Range("H138").Value = 3720.5256
$ws.Range("J138").Value = 5172.3335
$ws.Range("L138").Value = 15517.0005
$ws.Range("N138").Value = -25797.0005

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4800.52
$ws.Range("I32").Value = 4653.1787
$ws.Range("K32").Value = 4653.1787
$ws.Range("M32").Value = -4366.1787
$ws.Range("H45").Value = 78648.22
$ws.Range("I45").Value = 98947.05
$ws.Range("K45").Value = 98947.05
$ws.Range("M45").Value = -98570.05
$ws.Range("H61").Value = 5493.8
$ws.Range("I61").Value = 5309.635
$ws.Range("J61").Value = 6230.4614
$ws.Range("K61").Value = 5309.635
$ws.Range("L61").Value = 6230.4614
$ws.Range("M61").Value = -5097.635
$ws.Range("N61").Value = -6654.4614
$ws.Range("H74").Value = 1600.475
$ws.Range("I74").Value = 879.9583
$ws.Range("K74").Value = 879.9583
$ws.Range("M74").Value = -5.958300000000008
$ws.Range("H77").Value = 1600.475
$ws.Range("I77").Value = 879.9583
$ws.Range("K77").Value = 4399.7915
$ws.Range("M77").Value = -31.79150000000027
$ws.Range("H110").Value = 2259.3333
$ws.Range("J110").Value = 2790.9092
$ws.Range("L110").Value = 2790.9092
$ws.Range("N110").Value = -6880.9092
$ws.Range("H132").Value = 1516.5
$ws.Range("I132").Value = 875.0
$ws.Range("K132").Value = 2625.0
$ws.Range("M132").Value = -95.0
$ws.Range("H136").Value = 5493.8
$ws.Range("I136").Value = 5309.635
$ws.Range("J136").Value = 6230.4614
$ws.Range("K136").Value = 15928.905
$ws.Range("L136").Value = 18691.3842
$ws.Range("M136").Value = -13378.905
$ws.Range("N136").Value = -23791.3842

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3665.3872
$ws.Range("I20").Value = 2814.6316
$ws.Range("K20").Value = 2814.6316
$ws.Range("M20").Value = -2567.6316
$ws.Range("H86").Value = 4141.758
$ws.Range("I86").Value = 4685.92
$ws.Range("J86").Value = 2441.25
$ws.Range("K86").Value = 4685.92
$ws.Range("L86").Value = 2441.25
$ws.Range("M86").Value = -3562.92
$ws.Range("N86").Value = -4687.25
$ws.Range("H89").Value = 4141.758
$ws.Range("I89").Value = 4685.92
$ws.Range("J89").Value = 2441.25
$ws.Range("K89").Value = 23429.6
$ws.Range("L89").Value = 12206.25
$ws.Range("M89").Value = -17813.6
$ws.Range("N89").Value = -23438.25
$ws.Range("H107").Value = 938.8421
$ws.Range("I107").Value = 829.8571
$ws.Range("J107").Value = 1244.0
$ws.Range("K107").Value = 829.8571
$ws.Range("L107").Value = 1244.0
$ws.Range("M107").Value = 1090.1429
$ws.Range("N107").Value = -5084.0

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6590.0
$ws.Range("I31").Value = 6750.625
$ws.Range("J31").Value = 5733.3335
$ws.Range("K31").Value = 6750.625
$ws.Range("L31").Value = 5733.3335
$ws.Range("M31").Value = -6455.625
$ws.Range("N31").Value = -6323.3335
$ws.Range("H34").Value = 6590.0
$ws.Range("I34").Value = 6750.625
$ws.Range("J34").Value = 5733.3335
$ws.Range("K34").Value = 6750.625
$ws.Range("L34").Value = 5733.3335
$ws.Range("M34").Value = -6548.625
$ws.Range("N34").Value = -6137.3335
$ws.Range("H122").Value = 16688.5
$ws.Range("I122").Value = 18744.143
$ws.Range("K122").Value = 56232.429
$ws.Range("M122").Value = -53782.429
$ws.Range("H132").Value = 1561.0476
$ws.Range("J132").Value = 2869.8
$ws.Range("L132").Value = 8609.400000000001
$ws.Range("N132").Value = -13669.4
$ws.Range("H141").Value = 121955.766
$ws.Range("J141").Value = 125539.0
$ws.Range("L141").Value = 125539.0
$ws.Range("N141").Value = -135899.0

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 49627080.0
$ws.Range("I4").Value = 39559092.0
$ws.Range("K4").Value = 118677276.0
$ws.Range("M4").Value = -118677164.0
$ws.Range("H5").Value = 314401.28
$ws.Range("I5").Value = 879.4
$ws.Range("J5").Value = 456911.22
$ws.Range("K5").Value = 2638.2
$ws.Range("L5").Value = 1370733.66
$ws.Range("M5").Value = -2526.2
$ws.Range("N5").Value = -1370957.66
$ws.Range("H34").Value = 2213.5454
$ws.Range("J34").Value = 2915.3076
$ws.Range("L34").Value = 8745.9228
$ws.Range("N34").Value = -8913.9228
$ws.Range("H55").Value = 5809.353
$ws.Range("J55").Value = 10711.889
$ws.Range("L55").Value = 32135.667
$ws.Range("N55").Value = -32489.667
$ws.Range("H56").Value = 7136.057
$ws.Range("I56").Value = 7136.057
$ws.Range("K56").Value = 7136.057
$ws.Range("M56").Value = -6606.057
$ws.Range("H131").Value = 16396136.0
$ws.Range("I131").Value = 62505240.0
$ws.Range("K131").Value = 187515720.0
$ws.Range("M131").Value = -187510680.0
$ws.Range("H135").Value = 314401.28
$ws.Range("I135").Value = 879.4
$ws.Range("J135").Value = 456911.22
$ws.Range("K135").Value = 7914.599999999999
$ws.Range("L135").Value = 4112200.98
$ws.Range("M135").Value = -5379.599999999999
$ws.Range("N135").Value = -4117270.98

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10251.667
$ws.Range("I70").Value = 8017.875
$ws.Range("J70").Value = 14719.25
$ws.Range("K70").Value = 8017.875
$ws.Range("L70").Value = 14719.25
$ws.Range("M70").Value = -7747.875
$ws.Range("N70").Value = -15259.25
$ws.Range("H73").Value = 10251.667
$ws.Range("I73").Value = 8017.875
$ws.Range("J73").Value = 14719.25
$ws.Range("K73").Value = 8017.875
$ws.Range("L73").Value = 14719.25
$ws.Range("M73").Value = -7081.875
$ws.Range("N73").Value = -16591.25
$ws.Range("H132").Value = 2039.5555
$ws.Range("I132").Value = 2060.862
$ws.Range("J132").Value = 1951.2858
$ws.Range("K132").Value = 6182.586
$ws.Range("L132").Value = 5853.857400000001
$ws.Range("M132").Value = -3652.586
$ws.Range("N132").Value = -10913.8574

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1765.0625
$ws.Range("I16").Value = 1223.5172
$ws.Range("K16").Value = 1223.5172
$ws.Range("M16").Value = -1053.5172
$ws.Range("H55").Value = 874.4483
$ws.Range("I55").Value = 407.7647
$ws.Range("K55").Value = 407.7647
$ws.Range("M55").Value = -234.7647

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 31100.412
$ws.Range("I126").Value = 48680.0
$ws.Range("J126").Value = 5986.7144
$ws.Range("K126").Value = 146040.0
$ws.Range("L126").Value = 17960.1432
$ws.Range("M126").Value = -143570.0
$ws.Range("N126").Value = -22900.1432
